# Fruta / hortaliza, semanal
# Insert two new weekly rows (before what is currently row 516), shifting the
# existing "Brócoli" records down by two rows, then populate the two new
# rows with the latest week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows above row 516; everything below (old 516:635) shifts
# down to 518:637 and the sheet's used range grows to A1:R637.
$ws.Rows("516:517").Insert()

# Seed the two new rows with the same record shape (Mercado/Región/Fecha
# metadata columns) as the rows now sitting at 518:519, then overwrite the
# week-specific figures below.
$ws.Range("A518:R519").Copy()
$ws.Range("A516").PasteSpecial()

# Row 516 - "Primera" quality, new week
$ws.Range("D516").Value = 45275
$ws.Range("J516").Value = 1000
$ws.Range("K516").Value = 1000
$ws.Range("L516").Value = 1000
$ws.Range("M516").Value = 1000
$ws.Range("P516").Value = 1000

# Row 517 - "Segunda" quality, new week
$ws.Range("D517").Value = 45275
$ws.Range("J517").Value = 1000
$ws.Range("K517").Value = 800
$ws.Range("L517").Value = 800
$ws.Range("M517").Value = 800
$ws.Range("P517").Value = 800

$excel.CutCopyMode = 0
